# Importer: prevent product to be linked to variation child
# Insert a new test row (sheet row 46) with a product that tries to
# link a variation to a child product, pushing the previous rows 46-50
# down to 47-51.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift rows 46:50 down by inserting a new blank row at position 46.
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with the new test data.
$ws.Cells.Item(46, 2).Value = 42
$ws.Cells.Item(46, 3).Value = 28
$ws.Cells.Item(46, 4).Value = "This tries to link variation to child"
$ws.Cells.Item(46, 6).Value = "Color/Black"
$ws.Cells.Item(46, 7).Value = "Size/XS"
$ws.Cells.Item(46, 9).Value = 12
$ws.Cells.Item(46, 10).Value = 100
$ws.Cells.Item(46, 11).Value = "Test Category"
$ws.Cells.Item(46, 12).Value = "Test Category"
$ws.Cells.Item(46, 13).Value = "shirt1.jpeg"
$ws.Cells.Item(46, 14).Value = "shirt2.jpeg,shirt3.jpeg"

# Reselect the cell below the new data, matching the saved view state.
$ws.Range("B52").Select()
